$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-13)
# from serial 45174 (2023-09-05) to serial 45175 (2023-09-06)
$ws.Range("C2:C13").Value = 45175
